$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Chart title text: "Parser Comparison: Default ..." -> "Default ..."
# -----------------------------------------------------------------
$chart = $ws.ChartObjects(1).Chart
$chart.ChartTitle.Text = "Default VDM CSV, SAFE-CSV: Native and SAFE-CSV: Univocity"

# -----------------------------------------------------------------
# 2. Move/resize the chart (same size, new anchor position)
# -----------------------------------------------------------------
$co = $ws.ChartObjects(1)
$co.Left = 655.9661918676181
$co.Top = 34.941259842519685

# -----------------------------------------------------------------
# 3. Add the new "performance" block: headers (row 14), data (rows
#    15-24, mirroring A2:D11 but with C/D recomputed as ratios), and
#    a totals row (row 25) with AVERAGE formulas.
# -----------------------------------------------------------------
$ws.Range("A14").Value = "Count"
$ws.Range("B14").Value = "VDM CSV"
$ws.Range("C14").Value = "SAFE-CSV: Native"
$ws.Range("D14").Value = "SAFE-CSV: Univocity"

$counts = @(100, 500, 1000, 5000, 10000, 15000, 20000, 25000, 30000, 35000)
$vdm    = @(38, 154, 333, 4351, 15538, 33716, 62374, 94435, 135664, 186057)

for ($i = 0; $i -lt 10; $i++) {
    $r = 15 + $i
    $srcRow = 2 + $i
    $ws.Cells.Item($r, 1).Value = $counts[$i]
    $ws.Cells.Item($r, 2).Value = $vdm[$i]
    $ws.Cells.Item($r, 3).Formula = "=B$srcRow/C$srcRow"
    $ws.Cells.Item($r, 4).Formula = "=B$srcRow/D$srcRow"
}

# -----------------------------------------------------------------
# 4. Turn A14:D24 into a second table ("Table13") with a totals row.
# -----------------------------------------------------------------
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A14:D24"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table13"
$lo.ShowTotals = $true

$ws.Range("C25").Formula = "=AVERAGE(Table13[SAFE-CSV: Native])"
$ws.Range("D25").Formula = "=AVERAGE(Table13[SAFE-CSV: Univocity])"

# -----------------------------------------------------------------
# 5. View tweaks: zoom 145 -> 85, selection D15 -> D26, no frozen
#    top-left cell override.
# -----------------------------------------------------------------
$excel.ActiveWindow.Zoom = 85
[void]$ws.Range("D26").Select()
